$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A65").Value = "2025/12/05 11:00"
$ws.Range("B65").Value = "-"
$ws.Range("C65").Value = "-"
$ws.Range("D65").Value = "-"
$ws.Range("E65").Value = "-"
$ws.Range("F65").Value = "-"
$ws.Range("G65").Value = "-"
